# Append new ERA data rows (865-880) to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records: [rank, name, era]
$newRows = @(
    @(863, "ángel perdomo", 3.72),
    @(864, "brent honeywell", 4.82),
    @(865, "carl edwards", 3.69),
    @(866, "daniel lynch", 4.64),
    @(867, "duane underwood", 5.18),
    @(868, "j.b. bukauskas", 1.29),
    @(869, "j.t. chargois", 3.61),
    @(870, "jaime barría", 5.68),
    @(871, "jose ferrer", 5.03),
    @(872, "julio teherán", 4.4),
    @(873, "mark leiter", 3.5),
    @(874, "matt boyd", 5.45),
    @(875, "matthew bowman", 9),
    @(876, "mike king", 2.75),
    @(877, "néstor cortés", 4.97),
    @(878, "nick martínez", 3.43)
)

$startRow = 865
$lastExistingRow = 864

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rec = $newRows[$i]

    # Copy formatting (styles) from the last data row onto the new row first,
    # so column A keeps style index 1 (bordered/bold/centered) like the rest.
    $srcRange = $ws.Range("A" + $lastExistingRow + ":C" + $lastExistingRow)
    $dstRange = $ws.Range("A" + $rowIndex + ":C" + $rowIndex)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
}
